$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "length_pipe"
$ws.Range("A3").Value = "length_pipe"
$ws.Range("B2").Value = [double]"0.007549867130840628"
$ws.Range("C2").Value = [double]"1.299046696460312"
$ws.Range("D2").Value = [double]"0.8534736798886241"
$ws.Range("E2").Value = [double]"1.472183622241831"
$ws.Range("F2").Value = [double]"-12.24395991558717"
$ws.Range("G2").Value = [double]"3.468721592776412"
$ws.Range("K2").Value = [double]"0.0196"
$ws.Range("P2").Value = [double]"-0.3917021426743598"
$ws.Range("Q2").Value = [double]"-0.1039204065152347"
$ws.Range("B3").Value = [double]"3.519930706972955e-05"
$ws.Range("G3").Value = [double]"3.468721592776412"
$ws.Range("H3").Value = [double]"16.95283206967629"
$ws.Range("K3").Value = [double]"0.0196"
$ws.Range("B4").Value = [double]"3.519930706972955e-05"
$ws.Range("G4").Value = [double]"3.174724336146332"
$ws.Range("K4").Value = [double]"0.0196"
$ws.Range("B5").Value = [double]"3.519930706972955e-05"
$ws.Range("G5").Value = [double]"6.90979573"
$ws.Range("J5").Value = [double]"7.281418596835137"
$ws.Range("K5").Value = [double]"0.0196"
$ws.Range("B6").Value = [double]"3.042547539758451e-05"
$ws.Range("G6").Value = [double]"3.468721592776412"
$ws.Range("K6").Value = [double]"0.02267528804567081"
$ws.Range("B7").Value = [double]"3.519930706972955e-05"
$ws.Range("G7").Value = [double]"3.468721592776412"
$ws.Range("K7").Value = [double]"0.0196"
$ws.Range("B8").Value = [double]"4.679725849344545e-05"
$ws.Range("G8").Value = [double]"3.468721592776412"
$ws.Range("K8").Value = [double]"0.0196"
$ws.Range("B9").Value = [double]"4.25979166388954e-05"
$ws.Range("G9").Value = [double]"3.468721592776412"
$ws.Range("K9").Value = [double]"0.0196"
$ws.Range("B10").Value = [double]"3.329065653081932e-05"
$ws.Range("G10").Value = [double]"3.468721592776412"
$ws.Range("K10").Value = [double]"0.0196"
$ws.Range("B11").Value = [double]"3.310281177589487e-05"
$ws.Range("G11").Value = [double]"3.468721592776412"
$ws.Range("K11").Value = [double]"0.0196"
$ws.Range("B12").Value = [double]"3.068538462004525e-05"
$ws.Range("G12").Value = [double]"3.468721592776412"
$ws.Range("K12").Value = [double]"0.0196"
$ws.Range("B13").Value = [double]"3.300896183567589e-05"
$ws.Range("G13").Value = [double]"3.468721592776412"
$ws.Range("K13").Value = [double]"0.0196"
